$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws1.Range("B8").Value = "2023-09-01T08:48:57+00:00"

# Set the Case Sensitive value (row 14, column B) to the literal text "true".
# A direct Value assignment of "true" gets auto-coerced to a Boolean by
# Excel's type-sniffing, so instead stage the literal text as a formula
# result in a scratch cell (forces a text result), copy it, and paste
# only the value into the target cell so the destination keeps its
# original number format/style and becomes a real text cell.
$ws1.Range("Z1").Formula = '="true"'
$ws1.Range("Z1").Copy()
$ws1.Range("B14").PasteSpecial(-4163) # xlPasteValues
$ws1.Range("Z1").Clear()

# --- Concepts sheet ---
$ws2 = $wb.Worksheets.Item("Concepts")

# Add the Definition text for the phase-3-phase-4 concept (row 2, column D)
$ws2.Range("D2").Value = "Trials that are a combination of phases III and IV."
